$wb = $excel.ActiveWorkbook

# The new per-market sheets are created by copying the existing "Denmark"
# sheet (the last tab, which carries all the static boilerplate / merged
# cells / styles shared by every "XXX Market" test-data sheet) and then
# updating the two cells that vary per market: B4 (ticket reference) and
# B2 (market name).

function Add-MarketSheet($afterSheetName, $newName, $ticket, $marketName) {
    $src = $wb.Worksheets.Item($afterSheetName)
    $src.Copy($null, $src)
    $new = $wb.Worksheets.Item($wb.Worksheets.Count)
    $new.Name = $newName

    # Ticket reference first, then market name - matches the order the
    # corresponding shared strings were appended in the authored workbook.
    $new.Range("B4").Value = $ticket
    $new.Range("B2").Value = $marketName

    # Whole-sheet selection, as left behind in the saved file.
    $new.Range("A1:D18").Select()

    return $new
}

Add-MarketSheet "Denmark" "Russia" "NGC-2929/T2907/T2898" "Russia Market" | Out-Null
Add-MarketSheet "Russia" "Finland" "NGC-3130/T2884" "Finland Market" | Out-Null
$hungary = Add-MarketSheet "Finland" "Hungary" "NGC-3104/T2976/T2990" "Hungary Market"

# Hungary is now the last tab - make it the active sheet/tab, matching the
# original author ending their session there.
$hungary.Activate()
